$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 47: add a (blank, yellow-filled) status cell in column G
# ---------------------------------------------------------------------------
$ws.Range("G5").Copy()
$ws.Range("G47").PasteSpecial(-4122)
$ws.Range("G47").Value = ""

# ---------------------------------------------------------------------------
# Row 48: add date in column F and a (blank, yellow-filled) status cell in G
# ---------------------------------------------------------------------------
$ws.Range("F45").Copy()
$ws.Range("F48").PasteSpecial(-4122)
$ws.Range("F48").Value = 41926

$ws.Range("G5").Copy()
$ws.Range("G48").PasteSpecial(-4122)
$ws.Range("G48").Value = ""

# ---------------------------------------------------------------------------
# Row 49: add date in column F and a (blank, red-filled) status cell in G
# ---------------------------------------------------------------------------
$ws.Range("F45").Copy()
$ws.Range("F49").PasteSpecial(-4122)
$ws.Range("F49").Value = 41926

$ws.Range("G2").Copy()
$ws.Range("G49").PasteSpecial(-4122)
$ws.Range("G49").Value = ""

# ---------------------------------------------------------------------------
# Row 50 (new) - copy formats from similar existing rows, and set the
# non-text (date) values. Style indices used: A/F=1 (date), D/E=2 (wrap),
# G=5 (green), matching the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("A45").Copy()
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("A50").Value = 41928

$ws.Range("D45").Copy()
$ws.Range("D50").PasteSpecial(-4122)

$ws.Range("E45").Copy()
$ws.Range("E50").PasteSpecial(-4122)

$ws.Range("F45").Copy()
$ws.Range("F50").PasteSpecial(-4122)
$ws.Range("F50").Value = 41928

$ws.Range("G3").Copy()
$ws.Range("G50").PasteSpecial(-4122)
$ws.Range("G50").Value = ""

$ws.Rows.Item(50).RowHeight = 180

# ---------------------------------------------------------------------------
# Row 51 (new) - G=6 (yellow)
# ---------------------------------------------------------------------------
$ws.Range("A45").Copy()
$ws.Range("A51").PasteSpecial(-4122)
$ws.Range("A51").Value = 41928

$ws.Range("D45").Copy()
$ws.Range("D51").PasteSpecial(-4122)

$ws.Range("E45").Copy()
$ws.Range("E51").PasteSpecial(-4122)

$ws.Range("F45").Copy()
$ws.Range("F51").PasteSpecial(-4122)
$ws.Range("F51").Value = 41928

$ws.Range("G5").Copy()
$ws.Range("G51").PasteSpecial(-4122)
$ws.Range("G51").Value = ""

$ws.Rows.Item(51).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 52 (new) - G=5 (green)
# ---------------------------------------------------------------------------
$ws.Range("A45").Copy()
$ws.Range("A52").PasteSpecial(-4122)
$ws.Range("A52").Value = 41928

$ws.Range("D45").Copy()
$ws.Range("D52").PasteSpecial(-4122)

$ws.Range("E45").Copy()
$ws.Range("E52").PasteSpecial(-4122)

$ws.Range("F45").Copy()
$ws.Range("F52").PasteSpecial(-4122)
$ws.Range("F52").Value = 41928

$ws.Range("G3").Copy()
$ws.Range("G52").PasteSpecial(-4122)
$ws.Range("G52").Value = ""

$ws.Rows.Item(52).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 53 (new) - G=5 (green)
# ---------------------------------------------------------------------------
$ws.Range("A45").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$ws.Range("A53").Value = 41928

$ws.Range("D45").Copy()
$ws.Range("D53").PasteSpecial(-4122)

$ws.Range("E45").Copy()
$ws.Range("E53").PasteSpecial(-4122)

$ws.Range("F45").Copy()
$ws.Range("F53").PasteSpecial(-4122)
$ws.Range("F53").Value = 41928

$ws.Range("G3").Copy()
$ws.Range("G53").PasteSpecial(-4122)
$ws.Range("G53").Value = ""

$ws.Rows.Item(53).RowHeight = 45

# ---------------------------------------------------------------------------
# Now assign the new shared-string text values, in the exact order the
# original author entered them (this controls the resulting shared-string
# table order/index assignment so it matches the target workbook exactly).
# ---------------------------------------------------------------------------
$v158 = @"
Mohammad Abed Anwarzai <abedafg@hotmail.com>
"@
$ws.Range("C50").Value = $v158
$v159 = @"
Thank you very much , I have used the SAM beta version for solar annual  energy production in Afghanistan.
But the energy production was negative for months ( May >>>August). For more detail you can see the attached picture. I am going to try to learn more and optimize the model.
Note : The Afghanistan data is not exist in SAM library . I have created the CSV file from NREL TMY files, but I couldn’t understand regarding {Beam and Pwp}. If possible please could you include the Afghanistan data as well (request ) .

"@
$ws.Range("D50").Value = $v159
$v160 = @"
Seems like a weather data issue. Asked for copy of files
"@
$ws.Range("E50").Value = $v160
$ws.Range("B50").Value = "Email from SAM"
$v161 = @"
I’m a student at Brunel University and I plan to use SAM for my studies in MSc Renewable Energy Engineering, today I’ve tried to download and install the software to my computer however the SAM key code I’ve received is not working so I need some help to figure out this issue.
"@
$ws.Range("D51").Value = $v161
$v162 = @"
batuhan büyükbaş <buyukbas.batuhan@gmail.com>
"@
$ws.Range("C51").Value = $v162
$v163 = @"
Asked for screenshot
"@
$ws.Range("E51").Value = $v163
$ws.Range("B51").Value = "Email"
$v164 = @"
Couldn’t get weather data to load from Midway Airport, Chicago, IL
Did not pick up rates for ComEd

"@
$ws.Range("D52").Value = $v164
$v165 = @"
Brittany Placek <fritzo31@comcast.net>
"@
$ws.Range("C52").Value = $v165
$v166 = @"
Confirmed issue with solar prospector "Midway Airport, Chicago, IL" but works with today's version. Utility rate alias issue.
"@
$ws.Range("E52").Value = $v166
$ws.Range("B52").Value = "Email from SAM"
$v167 = @"
This version of the program is very interesting, and I like the new interface. Does your team know the expected timeline to move this from a Beta program to the official new version? 
"@
$ws.Range("D53").Value = $v167
$v168 = @"
Response to Beta announcement
"@
$ws.Range("B53").Value = $v168
$v169 = @"
Elizabeth Youngblood <yohohoblood@gmail.com>
"@
$ws.Range("C53").Value = $v169
$v170 = @"
Responded
"@
$ws.Range("E53").Value = $v170

# ---------------------------------------------------------------------------
# Update selection to match the new active cell
# ---------------------------------------------------------------------------
$ws.Range("G53").Select()
